$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2
$ws.Range("A2").Value = 3
$ws.Range("A3").Value = 4

$ws.Range("A4").Select()
